# edit.ps1 - applies the "456a3b4" gh-pages data refresh to 上海-漫展信息.xlsx
# Sheet order: 1=展览, 2=演出, 3=本地生活, 4=全部类型

$wb = $excel.ActiveWorkbook

# ---- Sheet 1 (展览) : update "想去人数" (F column) values ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 1074
$ws.Range("F6").Value = 3313
$ws.Range("F10").Value = 729
$ws.Range("F15").Value = 644
$ws.Range("F16").Value = 1648
$ws.Range("F17").Value = 1648
$ws.Range("F18").Value = 30
$ws.Range("F19").Value = 317
$ws.Range("F22").Value = 619
$ws.Range("F25").Value = 616
$ws.Range("F26").Value = 76592
$ws.Range("F27").Value = 76593
$ws.Range("F29").Value = 646
$ws.Range("F30").Value = 33174
$ws.Range("F31").Value = 33174
$ws.Range("F32").Value = 457
$ws.Range("F33").Value = 11
$ws.Range("F35").Value = 43
$ws.Range("F37").Value = 934
$ws.Range("F38").Value = 246
$ws.Range("F40").Value = 517
$ws.Range("F41").Value = 1158
$ws.Range("F42").Value = 5377
$ws.Range("F43").Value = 726
$ws.Range("F47").Value = 341
$ws.Range("F48").Value = 8
$ws.Range("F49").Value = 4
$ws.Range("F50").Value = 16
$ws.Range("F51").Value = 38
$ws.Range("F52").Value = 3

# ---- Sheet 2 (演出) : update "想去人数" (F column) values ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("F15").Value = 975
$ws.Range("F16").Value = 6
$ws.Range("F18").Value = 396
$ws.Range("F23").Value = 488
$ws.Range("F35").Value = 1402
$ws.Range("F37").Value = 6
$ws.Range("F38").Value = 96
$ws.Range("F39").Value = 96

# ---- Sheet 2 (演出): insert a new event row at row 46 ----
# Before: rows 46-48 hold "石川绫子"/"新海诚"/"变形金刚" events.
# After: a new "游马晃祐粉丝见面会" event is inserted as row 46 and the
# three existing events shift down to rows 47-49. The "A" column is a
# simple sequential index (A<n> = n-1) independent of the shift, so after
# the physical row insert we restore A46:A49 to the correct sequence.
$ws.Rows.Item(46).Insert()

$ws.Range("A46").Value = 45
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "2024-06-23"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "上海·游马晃祐粉丝见面会"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "宜昌路179号 万代南梦宫上海文化中心"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2024.06.23 13:00-06.23 20:00"
$ws.Range("F46").Value = 2
$ws.Range("G46").Value = 480
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "https://show.bilibili.com/platform/detail.html?id=84330"
$ws.Range("I46").NumberFormat = "@"
$ws.Range("I46").Value = "//i2.hdslb.com/bfs/openplatform/202404/G0vOG4EZ1713257811188.jpeg"

$ws.Range("A47").Value = 46
$ws.Range("A48").Value = 47
$ws.Range("A49").Value = 48

# ---- Sheet 3 (本地生活) : update "想去人数" (F column) values ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 539
$ws.Range("F6").Value = 558

# ---- Sheet 4 (全部类型) : update "想去人数" (F column) values ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 539
$ws.Range("F6").Value = 1074
$ws.Range("F11").Value = 3313
$ws.Range("F16").Value = 729
$ws.Range("F19").Value = 558
$ws.Range("F24").Value = 1648
$ws.Range("F25").Value = 1649
$ws.Range("F26").Value = 30
$ws.Range("F27").Value = 317
$ws.Range("F28").Value = 6
$ws.Range("F31").Value = 619
$ws.Range("F33").Value = 616
$ws.Range("F34").Value = 76595
$ws.Range("F35").Value = 646
$ws.Range("F36").Value = 33174
$ws.Range("F37").Value = 457
$ws.Range("F38").Value = 11
$ws.Range("F40").Value = 43
$ws.Range("F43").Value = 246
$ws.Range("F45").Value = 517
$ws.Range("F47").Value = 5377
$ws.Range("F49").Value = 96
$ws.Range("F52").Value = 8
